# Apply the "mirrored row numbering" re-mapping to Grid_Wire_Channel_Mapping.
#
# Column A (rows 3-82) and Column E (rows 3-66) both get a simple, strictly
# increasing index: A3=0, A4=1, ... A82=79   and   E3=0, E4=1, ... E66=63.
#
# Column B only changes in a handful of rows where the "Grids 1" channel
# numbering was re-mapped (its other rows already match the new sequence).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: rows 3..82 -> 0..79 -------------------------------------
for ($row = 3; $row -le 82; $row++) {
    $ws.Cells.Item($row, 1).Value2 = $row - 3
}

# --- Column E: rows 3..66 -> 0..63 -------------------------------------
for ($row = 3; $row -le 66; $row++) {
    $ws.Cells.Item($row, 5).Value2 = $row - 3
}

# --- Column B: explicit re-mapped values for the affected rows ---------
$bChanges = @{
    19 = 64
    20 = 65
    21 = 66
    22 = 67
    39 = 68
    40 = 69
    41 = 70
    42 = 71
    59 = 72
    60 = 73
    61 = 74
    62 = 75
    79 = 76
    80 = 77
    81 = 78
    82 = 79
}

foreach ($row in $bChanges.Keys) {
    $ws.Cells.Item($row, 2).Value2 = $bChanges[$row]
}

# --- Update the active selection to match the saved session (C79) ------
$ws.Range("C79").Select()
